# This script reproduces the per-row cryptocurrency price/volume refresh
# captured by the commit: it rewrites the Price (D) and Volume(1h) (E) text
# for each data row, and for rows 43-44 it also swaps which coin (B/C) is
# listed, matching the data source's new ranking order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Many of the new Price values look like numbers (e.g. "1.000", "30.590.57")
# but must stay plain text exactly as scraped, so force the cells to Text
# format before writing them, then restore the default "Normal" style so the
# workbook formatting is left exactly as it was.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '30.590.57'
$ws.Range("E2").Value = '  -0.81%  '

$ws.Range("D3").Value = '1.875.87'
$ws.Range("E3").Value = '  -0.97%  '

$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("E5").Value = '  +1.06%  '

$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.01%  '

$ws.Range("D7").Value = '0.4741'
$ws.Range("E7").Value = '  -0.90%  '

$ws.Range("D8").Value = '0.2900'
$ws.Range("E8").Value = '  -0.37%  '

$ws.Range("D9").Value = '0.06477'
$ws.Range("E9").Value = '  -1.56%  '

$ws.Range("D10").Value = '21.94'
$ws.Range("E10").Value = '  +2.46%  '

$ws.Range("D11").Value = '0.07733'
$ws.Range("E11").Value = '  -0.72%  '

$ws.Range("D12").Value = '0.7380'
$ws.Range("E12").Value = '  -0.38%  '

$ws.Range("D13").Value = '1.876.37'
$ws.Range("E13").Value = '  -0.99%  '

$ws.Range("D14").Value = '95.82'
$ws.Range("E14").Value = '  -1.41%  '

$ws.Range("D15").Value = '5.172'
$ws.Range("E15").Value = '  -0.41%  '

$ws.Range("E16").Value = '  -2.78%  '

$ws.Range("D17").Value = '30.571.77'
$ws.Range("E17").Value = '  -0.83%  '

$ws.Range("D18").Value = '13.22'
$ws.Range("E18").Value = '  -2.70%  '

$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  +0.10%  '

$ws.Range("D20").Value = '0.000007470'
$ws.Range("E20").Value = '  -2.14%  '

$ws.Range("D21").Value = '2.121.30'
$ws.Range("E21").Value = '  -1.64%  '

$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.02%  '

$ws.Range("D23").Value = '5.209'
$ws.Range("E23").Value = '  -2.11%  '

$ws.Range("D24").Value = '6.168'
$ws.Range("E24").Value = '  -1.33%  '

$ws.Range("D25").Value = '165.48'
$ws.Range("E25").Value = '  -0.63%  '

$ws.Range("E26").Value = '  -2.19%  '

$ws.Range("D27").Value = '18.74'
$ws.Range("E27").Value = '  -2.28%  '

$ws.Range("D28").Value = '1.903'
$ws.Range("E28").Value = '  -4.48%  '

$ws.Range("D29").Value = '0.09879'
$ws.Range("E29").Value = '  -1.62%  '

$ws.Range("E30").Value = '  -2.58%  '

$ws.Range("D31").Value = '1.508'
$ws.Range("E31").Value = '  -0.83%  '

$ws.Range("D32").Value = '4.247'
$ws.Range("E32").Value = '  -3.38%  '

$ws.Range("D33").Value = '4.088'
$ws.Range("E33").Value = '  -1.35%  '

$ws.Range("D34").Value = '0.04772'
$ws.Range("E34").Value = '  -0.47%  '

$ws.Range("D35").Value = '1.120'
$ws.Range("E35").Value = '  -1.26%  '

$ws.Range("D36").Value = '0.6937'
$ws.Range("E36").Value = '  -1.89%  '

$ws.Range("D37").Value = '2.719'

$ws.Range("D38").Value = '0.01848'
$ws.Range("E38").Value = '  -1.53%  '

$ws.Range("D39").Value = '2.760'
$ws.Range("E39").Value = '  -0.50%  '

$ws.Range("D40").Value = '6.238'
$ws.Range("E40").Value = '  -3.35%  '

$ws.Range("D41").Value = '73.22'
$ws.Range("E41").Value = '  +3.14%  '

$ws.Range("D42").Value = '1.975'
$ws.Range("E42").Value = '  +1.96%  '

$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '0.4165'
$ws.Range("E43").Value = '  -1.41%  '

$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.00%  '

$ws.Range("D45").Value = '0.8355'
$ws.Range("E45").Value = '  -1.58%  '

$ws.Range("D46").Value = '101.57'
$ws.Range("E46").Value = '  -1.12%  '

$ws.Range("D47").Value = '9.333'
$ws.Range("E47").Value = '  -1.89%  '

$ws.Range("D48").Value = '35.34'
$ws.Range("E48").Value = '  -0.11%  '

$ws.Range("D49").Value = '6.959'
$ws.Range("E49").Value = '  -3.01%  '

$ws.Range("D50").Value = '922.37'
$ws.Range("E50").Value = '  -2.70%  '

$ws.Range("D51").Value = '0.05667'
$ws.Range("E51").Value = '  +0.60%  '

# Restore the original (unstyled) cell style now that the text values are set.
$dataRange.Style = "Normal"
